$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text content updates (sharedStrings) ---
$ws.Range("B2").Value = "rohan"
$ws.Range("L2").Value = "Reading ,Drawing"
$ws.Range("B3").Value = "mini"
$ws.Range("H3").Value = "abcd"
$ws.Range("L3").Value = "Reading ,Writing"

# --- Formatting updates ---
# Pincode/Phone columns (I, K) get their font color set to explicit black,
# matching the data rows' existing format; the header row cells (I1/K1) pick
# up the same formatting as the data cells below them.
$ws.Range("I2:I3").Font.Color = 0
$ws.Range("K2:K3").Font.Color = 0
$ws.Range("I2").Copy()
$ws.Range("I1").PasteSpecial(-4122)
$ws.Range("K2").Copy()
$ws.Range("K1").PasteSpecial(-4122)

# Row heights for the three data rows grew slightly
$ws.Rows("1:3").RowHeight = 19.5
